$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing cell content (values + formatting left intact structurally;
# row heights are handled explicitly below).
$ws.Cells.Clear()

# Remove the two trailing rows that no longer exist in the final layout
# (shrinks the sheet from A1:C26 to A1:C24).
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(25).Delete()

# Write cell values row by row.
# Row 1
$ws.Range("B1").Value = 'Ementa atual:'
$ws.Range("C1").Value = 'Ementa modificada (dados modificados em vermelho):'

# Row 2
$ws.Range("B2").Value = 'LOM3015'
$ws.Range("C2").Value = 'LOM3015'

# Row 3
$ws.Range("A3").Value = 'Nome:'
$ws.Range("B3").Value = ' Termodinâmica de Materiais'
$ws.Range("C3").Value = ' Termodinâmica de Materiais'

# Row 4
$ws.Range("A4").Value = 'Name:'
$ws.Range("B4").Value = 'Thermodynamics of Materials'
$ws.Range("C4").Value = 'Thermodynamics of Materials'

# Row 5
$ws.Range("A5").Value = 'Créditos-aula:'
$ws.Range("B5").Value = '4'
$ws.Range("C5").Value = '4'

# Row 6
$ws.Range("A6").Value = 'Créditos-trabalho'
$ws.Range("B6").Value = '0'
$ws.Range("C6").Value = '0'

# Row 7
$ws.Range("A7").Value = 'Carga horária:'
$ws.Range("B7").Value = '60 h'
$ws.Range("C7").Value = '60 h'

# Row 8
$ws.Range("A8").Value = 'Ativação:'
$ws.Range("B8").Value = '01/01/2011'
$ws.Range("C8").Value = '01/01/2011'

# Row 9
$ws.Range("A9").Value = 'Semestre ideal:'
$ws.Range("B9").Value = 'EM-4'
$ws.Range("C9").Value = 'EM-4'

# Row 10
$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = '3577649 - Carlos Angelo Nunes'
$ws.Range("C10").Value = '3577649 - Carlos Angelo Nunes'

# Row 11
$ws.Range("A11").Value = 'Objectives:'

# Row 12
$ws.Range("A12").Value = 'Docentes responsáveis:'

# Row 13
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = '01/01/2011'
$ws.Range("C13").Value = '01/01/2011'

# Row 14
$ws.Range("A14").Value = 'Short syllabus:'

# Row 15
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '3577649 - Carlos Angelo Nunes'
$ws.Range("C15").Value = '3577649 - Carlos Angelo Nunes'

# Row 16
$ws.Range("A16").Value = 'Syllabus:'

# Row 17
$ws.Range("A17").Value = 'Avaliação:'

# Row 18
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '1176388 - Luiz Tadeu Fernandes Eleno'
$ws.Range("C18").Value = '1176388 - Luiz Tadeu Fernandes Eleno'

# Row 19
$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'O curso será ministrado na forma de aulas expositivas.'
$ws.Range("C19").Value = 'O curso será ministrado na forma de aulas expositivas.'

# Row 20
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF). A nota final será calculada através da expressão:NF=(P1+2*P2)/3'
$ws.Range("C20").Value = 'Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF). A nota final será calculada através da expressão:NF=(P1+2*P2)/3'

# Row 21
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'Para os alunos que obtiverem 3,0≤NF<5,0, será aplicada uma avaliação de recuperação (R), com pontuação de 0 a 10, que levará ao cálculo da média final (MF) através da seguinte expressão:MF=(NF+R)/2'
$ws.Range("C21").Value = 'Para os alunos que obtiverem 3,0≤NF<5,0, será aplicada uma avaliação de recuperação (R), com pontuação de 0 a 10, que levará ao cálculo da média final (MF) através da seguinte expressão:MF=(NF+R)/2'

# Row 22
$ws.Range("A22").Value = 'Requisitos:'

# Row 23
$ws.Range("B23").Value = 'LOB1004 -  Cálculo II  (Requisito fraco)
'
$ws.Range("C23").Value = 'LOB1004 -  Cálculo II  (Requisito fraco)
'

# Row 24
$ws.Range("B24").Value = 'LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)
'
$ws.Range("C24").Value = 'LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)
'

# Fix up row heights to match the target layout.
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).EntireRow.AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).EntireRow.AutoFit()
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(24).RowHeight = 30
